$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "beer-onto:BeerVariables"
$ws.Range("B18").Value = "Beer Variables"
$ws.Range("C18").Value = "Beer related variables"

$ws.Range("A19").Value = "beer-onto:BeerSubjects"
$ws.Range("B19").Value = "Beer Subjects"
$ws.Range("C19").Value = "Beer related subjects"

$ws.Range("A20").Value = "beer-onto:BeerBreweries"
$ws.Range("B20").Value = "Beer Breweries"
